# "Add files via upload" - refreshed poker person-info stats (Dec upload).
# Updates session/points totals for several players and records a new
# "01. Royal Flush" best-hand result for Mark.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Maisy) ---
$ws.Range("D2").Value = 122
$ws.Range("I2").Value = 446
$ws.Range("J2").Value = 3.66

# --- Row 3 (Mark) ---
$ws.Range("D3").Value = 143
$ws.Range("I3").Value = 518
$ws.Range("J3").Value = 3.62
$ws.Range("L3").Value = "01. Royal Flush"
$ws.Range("M3").Value = "Ace,10,Jack,Queen,King"
$ws.Range("N3").Value = 19.12

# --- Row 4 (Matt) ---
$ws.Range("D4").Value = 191
$ws.Range("I4").Value = 715
$ws.Range("J4").Value = 3.74

# --- Row 6 (Pepe) ---
$ws.Range("D6").Value = 38
$ws.Range("I6").Value = 139
$ws.Range("J6").Value = 3.66

# --- Row 7 (Prashant) ---
$ws.Range("D7").Value = 142
$ws.Range("I7").Value = 613

# --- Row 8 (Jon) ---
$ws.Range("D8").Value = 196
$ws.Range("I8").Value = 698
$ws.Range("J8").Value = 3.56

# --- Row 9 (Andy) ---
$ws.Range("D9").Value = 85
$ws.Range("F9").Value = 10
$ws.Range("H9").Value = 10
$ws.Range("I9").Value = 294
$ws.Range("J9").Value = 3.46
$ws.Range("K9").Value = 19.11

# --- Row 10 (Anthony) ---
$ws.Range("D10").Value = 199
$ws.Range("I10").Value = 821
$ws.Range("J10").Value = 4.13

# --- Row 11 (Bob) ---
$ws.Range("D11").Value = 121
$ws.Range("I11").Value = 480
$ws.Range("J11").Value = 3.97
